# Commit: "إضافة حدث جديد في Card6 by admin at 2026-01-20 16:04:44"
#
# Adds a new maintenance-log event row to the bottom of the "Card6"
# worksheet's table (row 13), following the same layout used by every
# other row in that table:
#   A  = card number
#   B-K = (left blank for this event, as is the case for most rows)
#   L  = Date
#   M  = Event
#   N  = Correction
#   O  = Serviced by

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

$newRow = 13

# Column A repeats the card number ("6") found in every other row of
# this table. Force text formatting first so the numeric-looking string
# is stored as text, matching the rest of the column.
$colA = $ws.Cells.Item($newRow, 1)
$colA.NumberFormat = "@"
$colA.Value = "6"

# Columns B through K have no data for this event (same as most rows).
for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item($newRow, $col).Value = ""
}

# Date / Event / Correction / Serviced by for the new entry.
$ws.Cells.Item($newRow, 12).Value = "20/1/2026"
$ws.Cells.Item($newRow, 13).Value = "زياره توكيل"
$ws.Cells.Item($newRow, 14).Value = "تم تغير سوفت كرد لbc"
$ws.Cells.Item($newRow, 15).Value = "م. احمد علي توكيل"
